# Workbook has a single sheet "in"; grab it explicitly (falls back to
# ActiveSheet, which already points at it).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# F11 held the "日语近义词" prompt. Insert the new "回答样式... Final output are"
# clause right before "in the following format", matching the authored edit.
$newText = "日语近义词/提供20个和主题内容相似的日语单词，提供例句和中文翻译，讲解语法，具体说明使用上的差别。回答样式:[单词][翻译][例句][语法和用法说明]。Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3`n"

$ws.Range("F11").Value = $newText

# Longer wrapped text needs a taller row (231.75 -> 330.75pt), as in the diff.
$ws.Range("F11").RowHeight = 330.75

# Leave the selection parked on F11 (was F12), matching the saved cursor
# position after editing that cell.
$ws.Range("F11").Select()
